$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 881.6
$ws.Range("I39").Value = 736
$ws.Range("J39").Value = 1100
$ws.Range("K39").Value = 2208
$ws.Range("L39").Value = 3300
$ws.Range("M39").Value = -1912
$ws.Range("N39").Value = -3892
$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138
$ws.Range("H92").Value = 1742.25
$ws.Range("I92").Value = 1773.1111
$ws.Range("J92").Value = 1649.6666
$ws.Range("K92").Value = 1773.1111
$ws.Range("L92").Value = 1649.6666
$ws.Range("M92").Value = -525.1111000000001
$ws.Range("N92").Value = -4145.6666
$ws.Range("H113").Value = 3618.375
$ws.Range("I113").Value = 3124.25
$ws.Range("K113").Value = 3124.25
$ws.Range("M113").Value = 129.75
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 74.666664
$ws.Range("I4").Value = 65.59999999999999
$ws.Range("K4").Value = 65.59999999999999
$ws.Range("M4").Value = 50.40000000000001
$ws.Range("H32").Value = 11447.135
$ws.Range("I32").Value = 9804.362999999999
$ws.Range("K32").Value = 9804.362999999999
$ws.Range("M32").Value = -9517.362999999999
$ws.Range("H45").Value = 3162.5
$ws.Range("I45").Value = 2996
$ws.Range("K45").Value = 2996
$ws.Range("M45").Value = -2619
$ws.Range("H110").Value = 2413.76
$ws.Range("I110").Value = 1243.0454
$ws.Range("K110").Value = 1243.0454
$ws.Range("M110").Value = 801.9546
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3390.8572
$ws.Range("I64").Value = 1822
$ws.Range("J64").Value = 5482.6665
$ws.Range("K64").Value = 1822
$ws.Range("L64").Value = 5482.6665
$ws.Range("M64").Value = -1597
$ws.Range("N64").Value = -5932.6665
$ws.Range("H67").Value = 3390.8572
$ws.Range("I67").Value = 1822
$ws.Range("J67").Value = 5482.6665
$ws.Range("K67").Value = 1822
$ws.Range("L67").Value = 5482.6665
$ws.Range("M67").Value = -1042
$ws.Range("N67").Value = -7042.6665
$ws.Range("H75").Value = 13273.333
$ws.Range("I75").Value = 13273.333
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 13273.333
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -12337.333
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 13273.333
$ws.Range("I78").Value = 13273.333
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 39819.999
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -35139.999
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 825.9167
$ws.Range("I80").Value = 438.33334
$ws.Range("J80").Value = 1213.5
$ws.Range("K80").Value = 438.33334
$ws.Range("L80").Value = 1213.5
$ws.Range("M80").Value = 559.66666
$ws.Range("N80").Value = -3209.5
$ws.Range("H83").Value = 825.9167
$ws.Range("I83").Value = 438.33334
$ws.Range("J83").Value = 1213.5
$ws.Range("K83").Value = 2191.6667
$ws.Range("L83").Value = 6067.5
$ws.Range("M83").Value = 2800.3333
$ws.Range("N83").Value = -16051.5
$ws.Range("H86").Value = 13125
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
$ws.Range("H89").Value = 13125
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2043.3334
$ws.Range("I31").Value = 2043.3334
$ws.Range("K31").Value = 2043.3334
$ws.Range("M31").Value = -1748.3334
$ws.Range("H34").Value = 2043.3334
$ws.Range("I34").Value = 2043.3334
$ws.Range("K34").Value = 2043.3334
$ws.Range("M34").Value = -1841.3334
$ws.Range("H62").Value = 7745.8335
$ws.Range("I62").Value = 8997
$ws.Range("J62").Value = 1490
$ws.Range("K62").Value = 8997
$ws.Range("L62").Value = 1490
$ws.Range("M62").Value = -8373
$ws.Range("N62").Value = -2738
$ws.Range("H65").Value = 7745.8335
$ws.Range("I65").Value = 8997
$ws.Range("J65").Value = 1490
$ws.Range("K65").Value = 44985
$ws.Range("L65").Value = 7450
$ws.Range("M65").Value = -41865
$ws.Range("N65").Value = -13690
$ws.Range("H68").Value = 72295
$ws.Range("J68").Value = 72295
$ws.Range("L68").Value = 72295
$ws.Range("N68").Value = -73793
$ws.Range("H71").Value = 72295
$ws.Range("J71").Value = 72295
$ws.Range("L71").Value = 216885
$ws.Range("N71").Value = -224373
$ws.Range("H74").Value = 70814
$ws.Range("J74").Value = 70814
$ws.Range("L74").Value = 70814
$ws.Range("N74").Value = -72562
$ws.Range("H77").Value = 70814
$ws.Range("J77").Value = 70814
$ws.Range("L77").Value = 212442
$ws.Range("N77").Value = -221178
$ws.Range("H134").Value = 2532.2222
$ws.Range("I134").Value = 1977.3077
$ws.Range("K134").Value = 5931.9231
$ws.Range("M134").Value = -3396.9231
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4483.423
$ws.Range("J109").Value = 4999.905
$ws.Range("L109").Value = 14999.715
$ws.Range("N109").Value = -17079.715
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2142975
$ws.Range("I11").Value = 2500112.5
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 2500112.5
$ws.Range("L11").Value = 150
$ws.Range("M11").Value = -2499973.5
$ws.Range("N11").Value = -428
$ws.Range("H102").Value = 6000
$ws.Range("I102").Value = 6000
$ws.Range("K102").Value = 6000
$ws.Range("M102").Value = -4378
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4828.3335
$ws.Range("I7").Value = 5093
$ws.Range("K7").Value = 5093
$ws.Range("M7").Value = -4981
$ws.Range("H40").Value = 2011.625
$ws.Range("I40").Value = 1584.8572
$ws.Range("K40").Value = 1584.8572
$ws.Range("M40").Value = -1448.8572
$ws.Range("H61").Value = 1950
$ws.Range("J61").Value = 1900
$ws.Range("L61").Value = 1900
$ws.Range("N61").Value = -2304
$ws.Range("H113").Value = 1950
$ws.Range("J113").Value = 1900
$ws.Range("L113").Value = 1900
$ws.Range("N113").Value = -6240
$ws.Range("H126").Value = 4828.3335
$ws.Range("I126").Value = 5093
$ws.Range("K126").Value = 15279
$ws.Range("M126").Value = -12809
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 985.375
$ws.Range("I113").Value = 1173.7778
$ws.Range("J113").Value = 743.1429000000001
$ws.Range("K113").Value = 3521.3334
$ws.Range("L113").Value = 2229.4287
$ws.Range("M113").Value = -1351.3334
$ws.Range("N113").Value = -6569.4287
$ws.Range("H136").Value = 1472.8462
$ws.Range("I136").Value = 1195.1818
$ws.Range("K136").Value = 3585.5454
$ws.Range("M136").Value = -1035.5454
